$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("87811004_1121_MX", 78,   "MXN", 7567.7),
    @("87811004_1121_BR", 38,   "BRL", 440.44),
    @("87811004_1121_CA", 287,  "CAD", 1332.8),
    @("87811004_1121_LL", 43,   "USD", 112),
    @("87811004_1121_BG", 6,    "BGN", 16.04),
    @("87811004_1121_HU", 900,  "HUF", 1915878),
    @("87811004_1121_PE", 21,   "PEN", 174.3),
    @("87811004_1121_RO", 494,  "RON", 9795.57),
    @("87811004_1121_EU", 746,  "EUR", 3309.31),
    @("87811004_1121_PL", 42,   "PLN", 401.94),
    @("87811004_1121_DK", 22,   "DKK", 440.16),
    @("87811004_1121_GB", 400,  "GBP", 1309.43),
    @("87811004_1121_SE", 36,   "SEK", 919.87),
    @("87811004_1121_JP", 23,   "JPY", 7546),
    @("87811004_1121_CO", 31,   "COP", 347830),
    @("87811004_1121_CL", 35,   "CLP", 73248),
    @("87811004_1121_CZ", 23,   "CZK", 1718.21),
    @("87811004_1121_US", 1415, "USD", 8030.4),
    @("87811004_1121_CH", 78,   "CHF", 387.56),
    @("87811004_1121_NO", 26,   "NOK", 655.2),
    @("87811004_1121_NZ", 34,   "NZD", 147.6),
    @("87811004_1121_AU", 307,  "AUD", 1617.46)
)

$ws.Range("E2:E23").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $entry = $data[$i]
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = [string]$entry[3]
}
